# This script reproduces the "Updated cryptos list" commit: per-row price
# (column D) and 1h volume-change percentage (column E) refreshes, plus a few
# rows whose coin/link (columns B/C) got reordered/swapped.
#
# NumberFormat is forced to Text ("@") before writing any value that looks
# numeric so Excel stores the exact original string (e.g. "309.43", or
# "0.0000104") instead of silently parsing it into a binary double and losing
# the precise decimal text / trailing zeros / thousands-dot grouping. The
# style is then reset to "Normal" so the cell keeps its original (default)
# formatting/style index - only the stored value changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeRef, $value) {
    $rng = $ws.Range($rangeRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" '42.465.17'
Set-TextValue "E2" '  -0.22%  '

Set-TextValue "D3" '2.284.30'
Set-TextValue "E3" '  -0.28%  '

Set-TextValue "E4" '  +0.11%  '

Set-TextValue "D5" '309.43'
Set-TextValue "E5" '  -4.17%  '

Set-TextValue "D6" '102.64'
Set-TextValue "E6" '  -0.94%  '

Set-TextValue "D7" '0.621'
Set-TextValue "E7" '  -1.28%  '

Set-TextValue "E8" '  -0.20%  '

Set-TextValue "D9" '0.600'
Set-TextValue "E9" '  -1.42%  '

Set-TextValue "D10" '38.37'
Set-TextValue "E10" '  -3.87%  '

Set-TextValue "E11" '  -1.01%  '

Set-TextValue "D12" '8.18'
Set-TextValue "E12" '  -2.52%  '

Set-TextValue "E13" '  +0.51%  '

Set-TextValue "D14" '0.965'
Set-TextValue "E14" '  -0.54%  '

Set-TextValue "D15" '15.12'
Set-TextValue "E15" '  -0.54%  '

Set-TextValue "D16" '2.630.48'
Set-TextValue "E16" '  -0.29%  '

Set-TextValue "D17" '2.289.60'
Set-TextValue "E17" '  +0.05%  '

Set-TextValue "D18" '42.429.89'
Set-TextValue "E18" '  -0.13%  '

Set-TextValue "D19" '7.27'
Set-TextValue "E19" '  -1.75%  '

Set-TextValue "B20" 'ShibaInu'
Set-TextValue "C20" 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue "D20" '0.0000104'
Set-TextValue "E20" '  -1.37%  '

Set-TextValue "B21" 'InternetComputer(DFINITY)'
Set-TextValue "C21" 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue "D21" '13.42'
Set-TextValue "E21" '  -0.21%  '

Set-TextValue "D22" '72.95'
Set-TextValue "E22" '  -0.51%  '

Set-TextValue "D23" '269.17'
Set-TextValue "E23" '  -0.17%  '

Set-TextValue "D24" '3.38'
Set-TextValue "E24" '  -6.36%  '

Set-TextValue "E25" '  -3.12%  '

Set-TextValue "D26" '1.00'
Set-TextValue "E26" '  -0.35%  '

Set-TextValue "D27" '10.68'
Set-TextValue "E27" '  -2.21%  '

Set-TextValue "B28" 'Filecoin'
Set-TextValue "C28" 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue "D28" '6.95'
Set-TextValue "E28" '  +12.70%  '

Set-TextValue "B29" 'Toncoin'
Set-TextValue "C29" 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue "D29" '2.28'
Set-TextValue "E29" '  -2.21%  '

Set-TextValue "D30" '22.24'
Set-TextValue "E30" '  -1.35%  '

Set-TextValue "D31" '35.62'
Set-TextValue "E31" '  -6.53%  '

Set-TextValue "D32" '164.18'
Set-TextValue "E32" '  -0.42%  '

Set-TextValue "D33" '0.0846'
Set-TextValue "E33" '  -3.88%  '

Set-TextValue "E34" '  -3.02%  '

Set-TextValue "D35" '2.53'
Set-TextValue "E35" '  +1.08%  '

Set-TextValue "E36" '  -3.51%  '

Set-TextValue "D37" '4.47'
Set-TextValue "E37" '  -3.38%  '

Set-TextValue "D38" '0.0344'
Set-TextValue "E38" '  -3.06%  '

Set-TextValue "E39" '  +0.54%  '

Set-TextValue "D40" '3.60'
Set-TextValue "E40" '  -2.92%  '

Set-TextValue "B41" 'BitcoinSV'
Set-TextValue "C41" 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
Set-TextValue "D41" '112.28'
Set-TextValue "E41" '  +19.61%  '

Set-TextValue "B42" 'ARBITRUM'
Set-TextValue "C42" 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue "D42" '1.55'
Set-TextValue "E42" '  +1.36%  '

Set-TextValue "D43" '69.57'
Set-TextValue "E43" '  +0.20%  '

Set-TextValue "E44" '  -0.26%  '

Set-TextValue "E45" '  -0.60%  '

Set-TextValue "D46" '12.02'
Set-TextValue "E46" '  -2.60%  '

Set-TextValue "D47" '1.707.85'
Set-TextValue "E47" '  +6.66%  '

Set-TextValue "D48" '109.81'
Set-TextValue "E48" '  -2.50%  '

Set-TextValue "D49" '77.24'
Set-TextValue "E49" '  -5.45%  '

Set-TextValue "B50" 'THORChain'
Set-TextValue "C50" 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-TextValue "D50" '5.14'
Set-TextValue "E50" '  -2.74%  '

Set-TextValue "B51" 'FraxShare'
Set-TextValue "C51" 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue "D51" '8.62'
Set-TextValue "E51" '  -3.59%  '

